$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("N6").Value = 1.53
$ws.Range("O6").Value = 2.4

# Row 9
$ws.Range("G9").Value = 2.27
$ws.Range("H9").Value = 2.7
$ws.Range("I9").Value = 3.7
$ws.Range("J9").Value = 1.1
$ws.Range("K9").Value = 5.9
$ws.Range("L9").Value = 1.38
$ws.Range("M9").Value = 2.8
$ws.Range("N9").Value = 2.12
$ws.Range("O9").Value = 1.65
$ws.Range("P9").Value = 1.44
$ws.Range("Q9").Value = 2.6
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 1.95
$ws.Range("T9").Value = 6.9
$ws.Range("U9").Value = 11
$ws.Range("V9").Value = 8.5
$ws.Range("W9").Value = 25
$ws.Range("X9").Value = 19.5
$ws.Range("Y9").Value = 29
$ws.Range("Z9").Value = 5.9
$ws.Range("AA9").Value = 5.3
$ws.Range("AB9").Value = 12.5
$ws.Range("AC9").Value = 60
$ws.Range("AD9").Value = 500
$ws.Range("AE9").Value = 9.75
$ws.Range("AF9").Value = 21
$ws.Range("AG9").Value = 11.75
$ws.Range("AH9").Value = 60
$ws.Range("AI9").Value = 37
$ws.Range("AJ9").Value = 40

# Row 10
$ws.Range("G10").Value = 1.87
$ws.Range("H10").Value = 3.35
$ws.Range("I10").Value = 4.05
$ws.Range("J10").Value = 1.08
$ws.Range("K10").Value = 6.5
$ws.Range("L10").Value = 1.37
$ws.Range("M10").Value = 2.85
$ws.Range("N10").Value = 2.1
$ws.Range("O10").Value = 1.65
$ws.Range("P10").Value = 1.44
$ws.Range("Q10").Value = 2.6
$ws.Range("R10").Value = 1.93
$ws.Range("S10").Value = 1.78
$ws.Range("T10").Value = 6.2
$ws.Range("U10").Value = 8
$ws.Range("V10").Value = 8.5
$ws.Range("W10").Value = 15.5
$ws.Range("X10").Value = 16
$ws.Range("Y10").Value = 32
$ws.Range("Z10").Value = 6.5
$ws.Range("AA10").Value = 6.4
$ws.Range("AB10").Value = 17
$ws.Range("AC10").Value = 90
$ws.Range("AD10").Value = 800
$ws.Range("AF10").Value = 21
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 65
$ws.Range("AI10").Value = 40
$ws.Range("AJ10").Value = 50

# Row 19
$ws.Range("G19").Value = 3.45
$ws.Range("H19").Value = 2.95
$ws.Range("I19").Value = 2.18
$ws.Range("L19").Value = 1.35
$ws.Range("M19").Value = 2.7
$ws.Range("N19").Value = 2.02
$ws.Range("O19").Value = 1.62
$ws.Range("P19").Value = 1.47
$ws.Range("Q19").Value = 2.32
$ws.Range("R19").Value = 1.75
$ws.Range("S19").Value = 1.85
$ws.Range("T19").Value = 9.5
$ws.Range("V19").Value = 11.5
$ws.Range("W19").Value = 50
$ws.Range("X19").Value = 32
$ws.Range("Y19").Value = 40
$ws.Range("Z19").Value = 7.8
$ws.Range("AA19").Value = 5.7
$ws.Range("AB19").Value = 14
$ws.Range("AC19").Value = 70
$ws.Range("AD19").Value = 600
$ws.Range("AE19").Value = 6.7
$ws.Range("AF19").Value = 10.25
$ws.Range("AG19").Value = 8.75
$ws.Range("AH19").Value = 22
$ws.Range("AI19").Value = 19
$ws.Range("AJ19").Value = 30

# Row 21
$ws.Range("J21").Value = 1.11
$ws.Range("K21").Value = 6.5

# Row 23
$ws.Range("G23").Value = 2.77
$ws.Range("H23").Value = 3.25
$ws.Range("T23").Value = 7.8
$ws.Range("U23").Value = 12
$ws.Range("V23").Value = 8.75
$ws.Range("W23").Value = 26
$ws.Range("AB23").Value = 11
$ws.Range("AI23").Value = 14.5

# Row 29
$ws.Range("G29").Value = 2.4
$ws.Range("K29").Value = 7
$ws.Range("L29").Value = 1.44
$ws.Range("M29").Value = 2.63
$ws.Range("N29").Value = 2.4
$ws.Range("O29").Value = 1.53
$ws.Range("P29").Value = 1.5
$ws.Range("Q29").Value = 2.5
$ws.Range("S29").Value = 1.73
$ws.Range("T29").Value = 6.5
$ws.Range("V29").Value = 10
$ws.Range("W29").Value = 23
$ws.Range("X29").Value = 23
$ws.Range("Y29").Value = 41
$ws.Range("Z29").Value = 7
$ws.Range("AA29").Value = 6
$ws.Range("AB29").Value = 17
$ws.Range("AC29").Value = 67
$ws.Range("AE29").Value = 7.5
$ws.Range("AF29").Value = 13
$ws.Range("AH29").Value = 34
$ws.Range("AI29").Value = 29
$ws.Range("AJ29").Value = 41

# Row 30
$ws.Range("G30").Value = 3
$ws.Range("I30").Value = 2.1
$ws.Range("U30").Value = 17
$ws.Range("AG30").Value = 9.5
$ws.Range("AH30").Value = 21
$ws.Range("AI30").Value = 17

# Row 33
$ws.Range("J33").Value = 1.02
$ws.Range("K33").Value = 12
$ws.Range("N33").Value = 1.44
$ws.Range("O33").Value = 2.63

# Row 34
$ws.Range("N34").Value = 1.4
$ws.Range("O34").Value = 2.75

# Row 35
$ws.Range("G35").Value = 3.1
$ws.Range("H35").Value = 3.5
$ws.Range("I35").Value = 2.05
$ws.Range("N35").Value = 1.83
$ws.Range("O35").Value = 1.98
$ws.Range("U35").Value = 17
$ws.Range("V35").Value = 12
$ws.Range("W35").Value = 34
$ws.Range("X35").Value = 26
$ws.Range("Z35").Value = 12
$ws.Range("AC35").Value = 41
$ws.Range("AD35").Value = 450
$ws.Range("AE35").Value = 8.5
$ws.Range("AF35").Value = 11
$ws.Range("AH35").Value = 19
$ws.Range("AI35").Value = 17
